$d = $word.ActiveDocument

# --- Objetivos paragraph: split the run at two missing-space seams with manual line breaks ---
$d.Content.Find.Execute(
    "ductilidade dosmateriais metálicos",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "ductilidade dos^lmateriais metálicos",
    2) | Out-Null

$d.Content.Find.Execute(
    "técnicas de análiseenvolvidas. Habilitar",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "técnicas de análise^lenvolvidas. Habilitar",
    2) | Out-Null

# --- Programa paragraph: separate the heading word "Programa" from the numbered list with two breaks ---
$d.Content.Find.Execute(
    "Programa1.Importância",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Programa^l^l1.Importância",
    2) | Out-Null

# --- Bibliografia paragraph: put each numbered reference on its own line ---
$d.Content.Find.Execute(
    "286p.2.A. Garcia",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "286p.^l2.A. Garcia",
    2) | Out-Null

$d.Content.Find.Execute(
    "384p.3.C.A. Sciammarella",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "384p.^l3.C.A. Sciammarella",
    2) | Out-Null

$d.Content.Find.Execute(
    "460p.4.R.W. Hertzberg",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "460p.^l4.R.W. Hertzberg",
    2) | Out-Null

$d.Content.Find.Execute(
    "786p.5.C. Suryanarayana",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "786p.^l5.C. Suryanarayana",
    2) | Out-Null

$d.Content.Find.Execute(
    "450p.6.N.E. Dowling",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "450p.^l6.N.E. Dowling",
    2) | Out-Null

$d.Content.Find.Execute(
    "946p.7.Y. Lee",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "946p.^l7.Y. Lee",
    2) | Out-Null

$d.Content.Find.Execute(
    "402p.8.R.W. Evans",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "402p.^l8.R.W. Evans",
    2) | Out-Null

$d.Content.Find.Execute(
    "115p.9.L.F.M. Silva",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "115p.^l9.L.F.M. Silva",
    2) | Out-Null

$d.Content.Find.Execute(
    "391p.10.G.E. Dieter",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "391p.^l10.G.E. Dieter",
    2) | Out-Null
